$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2 changes from "001" to "002". A leading apostrophe forces text so the
# numeric-looking string is not coerced into a number; ClearFormats() then
# drops the "quote prefix" style flag that Value-assignment left behind so
# the cell's style index is unchanged (same as every other cell in the row).
$ws.Range("J2").Value = "'002"
$ws.Range("J2").ClearFormats()

# Notice date / report date
$ws.Range("M2").Value = "2020-12-18 00:00:00"
$ws.Range("N2").Value = "2020-06-30 00:00:00"

# Updated figures
$ws.Range("O2").Value = 1269538905.01
$ws.Range("P2").Value = 66064983.31
$ws.Range("Q2").Value = 271129928.17

# Ratio cells that become blank: assigning a lone "'" produces an empty but
# still-present text cell (same representation as the sheet's other blank
# inlineStr cells, e.g. AH2) instead of deleting the cell outright.
$ws.Range("R2").Value = "'"
$ws.Range("R2").ClearFormats()

$ws.Range("S2").Value = 277767010.42

$ws.Range("T2").Value = "'"
$ws.Range("T2").ClearFormats()

$ws.Range("U2").Value = 30525103.83

$ws.Range("V2").Value = "'"
$ws.Range("V2").ClearFormats()

$ws.Range("W2").Value = 586812417.6
$ws.Range("X2").Value = 238319310.03

$ws.Range("Y2").Value = "'"
$ws.Range("Y2").ClearFormats()

$ws.Range("Z2").Value = "'"
$ws.Range("Z2").ClearFormats()

$ws.Range("AA2").Value = "'"
$ws.Range("AA2").ClearFormats()

$ws.Range("AB2").Value = 682726487.41

$ws.Range("AC2").Value = "'"
$ws.Range("AC2").ClearFormats()

$ws.Range("AD2").Value = "'"
$ws.Range("AD2").ClearFormats()

$ws.Range("AE2").Value = "'"
$ws.Range("AE2").ClearFormats()

$ws.Range("AF2").Value = 172.6721338572
$ws.Range("AG2").Value = 46.2224840282
